# Trabalho1-LabES.docx: atualiza a data de entrega de 04/05 para 12/05
# e move o bookmark automatico "_GoBack" do Word para o final do novo
# texto (comportamento padrao do Word ao editar um trecho do documento).

$d = $word.ActiveDocument

# 1) Troca o texto da data de entrega.
[void]$d.Content.Find.Execute("04/05", $true, $false, $false, $false, $false,
                               $true, 1, $false, "12/05", 2)

# 2) Localiza o fim do texto recem-inserido ("12/05").
$found = $d.Content
[void]$found.Find.Execute("12/05", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$editEnd = $found.End

# 3) Remove o bookmark "_GoBack" existente (ele ficava logo apos "SIGAA").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 4) Insere o bookmark "_GoBack" colapsado logo apos "12/05". O motor
#    COM tem um problema ao adicionar um bookmark colapsado exatamente na
#    ultima posicao de um paragrafo (imediatamente antes da marca de
#    paragrafo); para contornar isso, insere-se temporariamente um
#    caractere apos o texto, cria-se o bookmark entre "12/05" e esse
#    caractere (posicao deixa de ser a ultima do paragrafo) e depois o
#    caractere temporario e removido, mantendo o bookmark colapsado no
#    lugar correto.
$placeholder = $d.Range($editEnd, $editEnd)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($editEnd, $editEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tempChar = $d.Range($editEnd, $editEnd + 1)
$tempChar.Delete()
